$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws5 = $wb.Worksheets.Add($null, $lastSheet)
$ws5.Name = "Test5"
$ws5.Range("A2:A13").Value = 1
$r = $ws5.Range("A2:A13")
$r.Font.ThemeColor = 1
$r.HorizontalAlignment = -4152
Write-Host "Done"
